# The author's email handle "halit.reiss.1" had its trailing "1"
# re-formatted with an explicit Calibri font, which splits the single
# text run "halit.reiss.1" into two runs:
#   "halit.reiss."  (unchanged formatting)
#   "1"             (explicit Calibri font)

$d = $word.ActiveDocument

# Locate the target text in the document body and narrow the range
# down to just the trailing "1" character.
$rng = $d.Content
$found = $rng.Find.Execute("halit.reiss.1")

if ($found) {
    $rng.Start = $rng.End - 1

    # Applying a font name to just this character splits the run and
    # records the new font on the new run's rPr, matching the diff.
    $rng.Font.Name = "Calibri"
}
